$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.142.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -5.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.232.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.14%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.34'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.584'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -8.60%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.563'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0828'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -10.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -10.17%  '

$ws.Range("E14").Value = '  -1.22%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.576.68'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.01%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.867'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -12.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.249.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.100.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.49%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0968'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.58%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.18'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -13.30%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '238.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.63%  '

$ws.Range("E27").Value = '  -0.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.08%  '

$ws.Range("E29").Value = '  -2.04%  '

$ws.Range("E30").Value = '  -10.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -15.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.54%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.42'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.11%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0879'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '152.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.68%  '

$ws.Range("E36").Value = '  -3.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.19'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.22%  '

$ws.Range("E39").Value = '  -7.91%  '

$ws.Range("E40").Value = '  -6.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.104'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -11.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0326'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.51%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.740.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.80%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.206'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.18%  '

$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.75'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -11.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.33'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -10.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.84'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -17.17%  '
